$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 376
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 20

$ws.Range("B5").Value = 0.9399999999999999
$ws.Range("C5").Value = 0.01
$ws.Range("D5").Value = 0.05
